# 글로벌 상수 TimeSecToGetOneTicket 960 추가
# Insert a new constant row "TimeSecToGetOneTicket" = 960 right after the
# existing "TimeSecToGetOneEnergy" row (row 5) on the GlobalConstantIntTable
# sheet, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "TimeSecToGetOneTicket"
$ws.Range("B6").Value = 960
